$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 296.5995795127912
$ws.Range("G2").Value = 23853195.83998579
$ws.Range("J2").Value = 155.5555284265615
$ws.Range("L2").Value = 29.77400213642043
$ws.Range("M2").Value = 56.7704446341448

$ws.Range("F4").Value = 254.6468579941883
$ws.Range("G4").Value = 86396806.08009991
$ws.Range("J4").Value = 136.4318069672803
$ws.Range("L4").Value = 26.22536219704341
$ws.Range("M4").Value = 48.94904078224418

$ws.Range("F5").Value = 109.5570717781737
$ws.Range("G5").Value = 10516814.6157781
$ws.Range("J5").Value = 70.94403448289002
$ws.Range("L5").Value = 12.33809295354609
$ws.Range("M5").Value = 19.05340378750846

$ws.Range("F6").Value = 188.1561850756708
$ws.Range("G6").Value = 15513484.0253405
$ws.Range("J6").Value = 87.29761195523272
$ws.Range("L6").Value = 15.39399160768768
$ws.Range("M6").Value = 33.17931234447459

$ws.Range("F9").Value = 152.0431323330019
$ws.Range("G9").Value = 3835150.174225295
$ws.Range("J9").Value = 76.76441501651911
$ws.Range("L9").Value = 11.96328545711986
$ws.Range("M9").Value = 23.69503361033796

$ws.Range("F10").Value = 203.8725298381634
$ws.Range("G10").Value = 17507776.58339982
$ws.Range("J10").Value = 109.6154306498862
$ws.Range("L10").Value = 17.30769957629782
$ws.Range("M10").Value = 32.19039944813106

$ws.Range("F12").Value = 394.2001942322248
$ws.Range("G12").Value = 109199621.6541033
$ws.Range("J12").Value = 237.8346121433357
$ws.Range("L12").Value = 36.65518800585989
$ws.Range("M12").Value = 60.754328822503

$ws.Range("F14").Value = 261.6614664607744
$ws.Range("G14").Value = 24357207.70323403
$ws.Range("J14").Value = 137.7047020761761
$ws.Range("L14").Value = 21.90850140547216
$ws.Range("M14").Value = 41.6297375418785

$ws.Range("F16").Value = 452.8546878821612
$ws.Range("G16").Value = 23784104.6119052
$ws.Range("J16").Value = 155.1474534370854
$ws.Range("L16").Value = 22.39197270907106
$ws.Range("M16").Value = 65.35917662576583

$ws.Range("F17").Value = 411.5681105642029
$ws.Range("G17").Value = 27886284.82320945
$ws.Range("J17").Value = 186.3632919203486
$ws.Range("L17").Value = 27.11727099479407
$ws.Range("M17").Value = 59.88627841879335

$ws.Range("F18").Value = 309.4825811310093
$ws.Range("G18").Value = 49244660.93664289
$ws.Range("J18").Value = 145.6360948634211
$ws.Range("L18").Value = 20.83346290277071
$ws.Range("M18").Value = 44.27194974634016

$ws.Range("F19").Value = 511.8247641200482
$ws.Range("G19").Value = 33247887.26297089
$ws.Range("J19").Value = 216.8811954531695
$ws.Range("L19").Value = 30.87828208531443
$ws.Range("M19").Value = 72.87063044689339

$ws.Range("F22").Value = 382.5674025459423
$ws.Range("G22").Value = 151358304.6767375
$ws.Range("J22").Value = 161.9513826655055
$ws.Range("L22").Value = 21.83614148298951
$ws.Range("M22").Value = 51.58212169158772

$ws.Range("F23").Value = 319.4937778679435
$ws.Range("G23").Value = 22962171.73611701
$ws.Range("J23").Value = 112.0149651504303
$ws.Range("L23").Value = 15.40037587950109
$ws.Range("M23").Value = 43.92559747459114

$ws.Range("F24").Value = 215.5050416447036
$ws.Range("G24").Value = 6168599.512652317
$ws.Range("J24").Value = 158.1570523460328
$ws.Range("L24").Value = 21.56687077445902
$ws.Range("M24").Value = 29.38705113336867

$ws.Range("F25").Value = 411.9436297660876
$ws.Range("G25").Value = 92571022.36643706
$ws.Range("J25").Value = 204.0393534962862
$ws.Range("L25").Value = 27.64307840844873
$ws.Range("M25").Value = 55.80977327342983

$ws.Range("F26").Value = 385.7288799402649
$ws.Range("G26").Value = 35545606.28596731
$ws.Range("J26").Value = 221.9381011861096
$ws.Range("L26").Value = 28.63717434659479
$ws.Range("M26").Value = 49.77146837938902

$ws.Range("F27").Value = 490.9815975211757
$ws.Range("G27").Value = 116266901.3094148
$ws.Range("J27").Value = 247.0651887609964
$ws.Range("L27").Value = 31.46952924456126
$ws.Range("M27").Value = 62.53798772388307

$ws.Range("F29").Value = 409.0416064941321
$ws.Range("G29").Value = 199164655.1709391
$ws.Range("J29").Value = 162.591355634512
$ws.Range("L29").Value = 20.92633984192375
$ws.Range("M29").Value = 52.64574880736

$ws.Range("F30").Value = 240.9175228078685
$ws.Range("G30").Value = 43019836.26603248
$ws.Range("J30").Value = 133.4200771805907
$ws.Range("L30").Value = 16.00258255407172
$ws.Range("M30").Value = 28.89596999885579

$ws.Range("F32").Value = 371.2418378297348
$ws.Range("G32").Value = 54566501.55957337
$ws.Range("J32").Value = 153.0145020037951
$ws.Range("L32").Value = 17.73211069551741
$ws.Range("M32").Value = 43.02142134894449

$ws.Range("F33").Value = 271.4578161450656
$ws.Range("G33").Value = 7405305.943708388
$ws.Range("J33").Value = 70.69409066243966
$ws.Range("L33").Value = 7.469683007134808
$ws.Range("M33").Value = 28.68279112740688

$ws.Range("F35").Value = 330.6426781546463
$ws.Range("G35").Value = 8798651.309854053
$ws.Range("J35").Value = 121.7974987521325
$ws.Range("L35").Value = 13.61605451463945
$ws.Range("M35").Value = 36.96339232533892
